# Draft profiles from the defunct "develop" branch
#
# Updates a handful of mapping notes on the "Data" sheet of the
# OutcomeOfCare - STU3 mapping workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Row 3: top-level OutcomeOfCare summary row -----------------------
$ws.Range("P3").Value = "CarePlan  / DiagnosticResult"
$ws.Range("R3").Value = "** OutcomeOfCare`n"
$ws.Rows.Item(3).RowHeight = 25.5

# --- Row 4: HealthcareResult row ---------------------------------------
$ws.Range("P4").Value = "CarePlan.activity:nursingIntervention.outcomeCodeableConcept Or derived profile on zib-TextResult."
$ws.Range("Q4").Value = ""

# --- Row 5: MeasurementValue::GeneralMeasurement row --------------------
$ws.Rows.Item(5).RowHeight = 89.25

# --- Row 7: Intervention::NursingIntervention row ------------------------
$ws.Range("P7").Value = "Careplan.activity / DiagnosticReport.extention.partOf"
$ws.Range("R7").Value = "Maybe not the most suitable extension."

# --- Sheet view / selection state ---------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1
$ws.Range("R13").Select() | Out-Null
